# Updates the cryptocurrency price/volume figures in the active
# workbook's active worksheet to reflect the latest scraped values
# (GitHub Actions scheduled refresh). Column D holds the Price text,
# column E holds the Volume(1h) percentage text (padded with two
# leading/trailing spaces). A handful of new Price values parse as
# plain decimals (one dot, no thousands separators), so without help
# Excel would coerce them to numbers on assignment; those cells are
# pre-formatted as Text ("@") so they keep storing the scraped string
# verbatim, matching every other Price cell in the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D14", "D20", "D21", "D23", "D24", "D25", "D29", "D30", "D32", "D38", "D40", "D44", "D47", "D48", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.070.07"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "3.116.41"
$ws.Range("E3").Value = "  +2.83%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "580.20"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "173.79"
$ws.Range("E6").Value = "  +3.17%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.112.50"
$ws.Range("E8").Value = "  +2.83%  "
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  -3.43%  "
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "37.25"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D16").Value = "3.631.66"
$ws.Range("E16").Value = "  +2.78%  "
$ws.Range("D17").Value = "67.043.56"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "3.117.58"
$ws.Range("E19").Value = "  +2.74%  "
$ws.Range("D20").Value = "16.18"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").Value = "487.71"
$ws.Range("E21").Value = "  +3.47%  "
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("D23").Value = "7.61"
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").Value = "84.44"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").Value = "13.35"
$ws.Range("E25").Value = "  +3.62%  "
$ws.Range("E26").Value = "  +4.26%  "
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "8.04"
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("D30").Value = "2.41"
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("D32").Value = "28.91"
$ws.Range("E32").Value = "  +3.15%  "
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "47.55"
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("E39").Value = "  +3.16%  "
$ws.Range("D40").Value = "50.16"
$ws.Range("E40").Value = "  +1.05%  "
$ws.Range("E41").Value = "  +2.25%  "
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").Value = "2.81"
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D45").Value = "2.846.66"
$ws.Range("E45").Value = "  +4.47%  "
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").Value = "384.71"
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("D48").Value = "137.06"
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("D50").Value = "25.19"
$ws.Range("E50").Value = "  +2.47%  "
$ws.Range("D51").Value = "2.23"
$ws.Range("E51").Value = "  +0.35%  "
